# Adds a new "2022-Q4" sheet (fund holdings detail) right after the
# "总计" summary sheet, and inserts a corresponding "2022-Q4" row at the
# top of the "总计" sheet's table (shifting the other quarters down).

function Set-TextValue($range, $val) {
    # Force a numeric-looking string to be stored as TEXT (not Number),
    # without leaving a lingering NumberFormat/quotePrefix style on the
    # cell (ClearFormats drops the style index back to default == 0).
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for 2022-Q4, shift the rest down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

# Re-apply the existing "index column" style (bold/bordered) from A3
# (which already carries it) onto the freshly inserted A2.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 14
$total.Range("D2").Value = 1.68

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" detail sheet right before "2022-Q3".
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Pull the header/index-column style (bold font + border, centered) from
# a cell that already has it, and stamp it onto the header row and the
# whole index column in one shot.
$wb.Worksheets.Item(1).Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "161914", "万家创业板2年定期开放混合A",      "8.14",  "99.81", "6.70", "0.5454", 8),
    @(1,  "009837", "华夏磐锐一年定期开放混合A",        "14.15", "75.21", "3.19", "0.4514", 8),
    @(2,  "003713", "英大睿盛灵活配置混合A",            "2.39",  "93.29", "7.01", "0.1675", 6),
    @(3,  "003714", "英大睿盛灵活配置混合C",            "2.39",  "93.29", "7.01", "0.1675", 6),
    @(4,  "010676", "光大保德信新机遇混合",              "2.86",  "90.93", "5.01", "0.1433", 7),
    @(5,  "161915", "万家创业板2年定期开放混合C",      "1.59",  "99.81", "6.70", "0.1065", 8),
    @(6,  "001607", "英大策略优选混合A",                "0.59",  "93.12", "5.69", "0.0336", 5),
    @(7,  "010540", "浙商智多金稳健一年持有期混合C",  "1.37",  "25.01", "1.05", "0.0144", 10),
    @(8,  "010539", "浙商智多金稳健一年持有期混合A",  "1.27",  "25.01", "1.05", "0.0133", 10),
    @(9,  "012522", "英大稳固增强核心一年持有混合C",  "1.05",  "23.17", "1.24", "0.0130", 8),
    @(10, "009838", "华夏磐锐一年定期开放混合C",        "0.39",  "75.21", "3.19", "0.0124", 8),
    @(11, "012521", "英大稳固增强核心一年持有混合A",  "0.63",  "23.17", "1.24", "0.0078", 8),
    @(12, "001899", "东海中证社会发展安全产业主题指数", "0.20", "93.35", "1.35", "0.0027", 10),
    @(13, "001608", "英大策略优选混合C",                "0.02",  "93.12", "5.69", "0.0011", 5)
)

$r = 2
foreach ($row in $rows) {
    $q4.Range("A$r").Value = $row[0]
    Set-TextValue $q4.Range("B$r") $row[1]
    $q4.Range("C$r").Value = $row[2]
    Set-TextValue $q4.Range("D$r") $row[3]
    Set-TextValue $q4.Range("E$r") $row[4]
    Set-TextValue $q4.Range("F$r") $row[5]
    Set-TextValue $q4.Range("G$r") $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r++
}

Write-Output "done"
